# Rename the "TEST####" Order ID values in column D to "TST####"
# (commit: "added Insert new CHEMProduct regression test" -
#  the Katalon regression input workbook's Order IDs were renamed from the
#  TEST#### naming convention to TST####).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $v = $cell.Value2
    if ($v -ne $null -and $v -like "TEST*") {
        $cell.Value = "TST" + $v.Substring(4)
    }
}
